# Natmi output update (per Dr Hou's advice): the Sending/Target cluster set grows
# from {FAPs, sCs} / {ECs, FAPs, sCs} to a full {ECs, FAPs, sCs} x {ECs, FAPs, sCs}
# grid (9 rows instead of 6), and every row's computed statistics are refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ucn2"
$ws.Cells.Item(2,3).Value = "Il10rb"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.3323133333333333
$ws.Cells.Item(2,8).Value = 0.99694
$ws.Cells.Item(2,9).Value = 0.2302327407427936
$ws.Cells.Item(2,10).Value = 0.2302327407427935
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 86.24110266666668
$ws.Cells.Item(2,14).Value = 258.723308
$ws.Cells.Item(2,15).Value = 0.861191559324194
$ws.Cells.Item(2,16).Value = 0.861191559324194
$ws.Cells.Item(2,17).Value = 28.65906829750223
$ws.Cells.Item(2,18).Value = 257.93161467752
$ws.Cells.Item(2,19).Value = 0.1982744930077693
$ws.Cells.Item(2,20).Value = 0.1982744930077693

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ucn2"
$ws.Cells.Item(3,3).Value = "Il10rb"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.3323133333333333
$ws.Cells.Item(3,8).Value = 0.99694
$ws.Cells.Item(3,9).Value = 0.2302327407427936
$ws.Cells.Item(3,10).Value = 0.2302327407427935
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 9.389419666666667
$ws.Cells.Item(3,14).Value = 28.168259
$ws.Cells.Item(3,15).Value = 0.09376142829643226
$ws.Cells.Item(3,16).Value = 0.09376142829643226
$ws.Cells.Item(3,17).Value = 3.120229347495556
$ws.Cells.Item(3,18).Value = 28.08206412746
$ws.Cells.Item(3,19).Value = 0.02158695061264652
$ws.Cells.Item(3,20).Value = 0.02158695061264651

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ucn2"
$ws.Cells.Item(4,3).Value = "Il10rb"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.3323133333333333
$ws.Cells.Item(4,8).Value = 0.99694
$ws.Cells.Item(4,9).Value = 0.2302327407427936
$ws.Cells.Item(4,10).Value = 0.2302327407427935
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 4.51108
$ws.Cells.Item(4,14).Value = 13.53324
$ws.Cells.Item(4,15).Value = 0.04504701237937385
$ws.Cells.Item(4,16).Value = 0.04504701237937385
$ws.Cells.Item(4,17).Value = 1.499092031733333
$ws.Cells.Item(4,18).Value = 13.4918282856
$ws.Cells.Item(4,19).Value = 0.01037129712237779
$ws.Cells.Item(4,20).Value = 0.01037129712237779

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ucn2"
$ws.Cells.Item(5,3).Value = "Il10rb"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.476438
$ws.Cells.Item(5,8).Value = 1.429314
$ws.Cells.Item(5,9).Value = 0.3300849395169671
$ws.Cells.Item(5,10).Value = 0.3300849395169671
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 86.24110266666668
$ws.Cells.Item(5,14).Value = 258.723308
$ws.Cells.Item(5,15).Value = 0.861191559324194
$ws.Cells.Item(5,16).Value = 0.861191559324194
$ws.Cells.Item(5,17).Value = 41.08853847230134
$ws.Cells.Item(5,18).Value = 369.7968462507121
$ws.Cells.Item(5,19).Value = 0.2842663637720492
$ws.Cells.Item(5,20).Value = 0.2842663637720492

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Ucn2"
$ws.Cells.Item(6,3).Value = "Il10rb"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.476438
$ws.Cells.Item(6,8).Value = 1.429314
$ws.Cells.Item(6,9).Value = 0.3300849395169671
$ws.Cells.Item(6,10).Value = 0.3300849395169671
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 9.389419666666667
$ws.Cells.Item(6,14).Value = 28.168259
$ws.Cells.Item(6,15).Value = 0.09376142829643226
$ws.Cells.Item(6,16).Value = 0.09376142829643226
$ws.Cells.Item(6,17).Value = 4.473476327147333
$ws.Cells.Item(6,18).Value = 40.261286944326
$ws.Cells.Item(6,19).Value = 0.03094923538825229
$ws.Cells.Item(6,20).Value = 0.03094923538825229

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ucn2"
$ws.Cells.Item(7,3).Value = "Il10rb"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.476438
$ws.Cells.Item(7,8).Value = 1.429314
$ws.Cells.Item(7,9).Value = 0.3300849395169671
$ws.Cells.Item(7,10).Value = 0.3300849395169671
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 4.51108
$ws.Cells.Item(7,14).Value = 13.53324
$ws.Cells.Item(7,15).Value = 0.04504701237937385
$ws.Cells.Item(7,16).Value = 0.04504701237937385
$ws.Cells.Item(7,17).Value = 2.14924993304
$ws.Cells.Item(7,18).Value = 19.34324939736
$ws.Cells.Item(7,19).Value = 0.01486934035666568
$ws.Cells.Item(7,20).Value = 0.01486934035666568

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Ucn2"
$ws.Cells.Item(8,3).Value = "Il10rb"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.6346286666666666
$ws.Cells.Item(8,8).Value = 1.903886
$ws.Cells.Item(8,9).Value = 0.4396823197402394
$ws.Cells.Item(8,10).Value = 0.4396823197402393
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 86.24110266666668
$ws.Cells.Item(8,14).Value = 258.723308
$ws.Cells.Item(8,15).Value = 0.861191559324194
$ws.Cells.Item(8,16).Value = 0.861191559324194
$ws.Cells.Item(8,17).Value = 54.73107599720978
$ws.Cells.Item(8,18).Value = 492.579683974888
$ws.Cells.Item(8,19).Value = 0.3786507025443756
$ws.Cells.Item(8,20).Value = 0.3786507025443756

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Ucn2"
$ws.Cells.Item(9,3).Value = "Il10rb"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.6346286666666666
$ws.Cells.Item(9,8).Value = 1.903886
$ws.Cells.Item(9,9).Value = 0.4396823197402394
$ws.Cells.Item(9,10).Value = 0.4396823197402393
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 9.389419666666667
$ws.Cells.Item(9,14).Value = 28.168259
$ws.Cells.Item(9,15).Value = 0.09376142829643226
$ws.Cells.Item(9,16).Value = 0.09376142829643226
$ws.Cells.Item(9,17).Value = 5.958794883830445
$ws.Cells.Item(9,18).Value = 53.629153954474
$ws.Cells.Item(9,19).Value = 0.04122524229553345
$ws.Cells.Item(9,20).Value = 0.04122524229553345

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Ucn2"
$ws.Cells.Item(10,3).Value = "Il10rb"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.6346286666666666
$ws.Cells.Item(10,8).Value = 1.903886
$ws.Cells.Item(10,9).Value = 0.4396823197402394
$ws.Cells.Item(10,10).Value = 0.4396823197402393
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 4.51108
$ws.Cells.Item(10,14).Value = 13.53324
$ws.Cells.Item(10,15).Value = 0.04504701237937385
$ws.Cells.Item(10,16).Value = 0.04504701237937385
$ws.Cells.Item(10,17).Value = 2.862860685626667
$ws.Cells.Item(10,18).Value = 25.76574617064
$ws.Cells.Item(10,19).Value = 0.01980637490033037
$ws.Cells.Item(10,20).Value = 0.01980637490033037

